# Applies the bilibili-scraped con-listing refresh: the expired
# 2024.02.24 Jingdezhen event (old row 2) drops off, every remaining
# event shifts up one row with refreshed interest counts, and the
# newly-uncovered last row is removed. Mirrors the same edit on both
# the "展览" and "全部类型" sheets (sheet1/sheet4), which are identical.
$wb = $excel.ActiveWorkbook

$rowsAfter = @{
  2 = @{
    B = "2024.03.02"
    C = "南昌·meeting动漫游戏嘉年华"
    D = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
    E = "2024.03.02 09:00-03.03 17:00"
    F = 1464
    G = 60
    H = "https://show.bilibili.com/platform/detail.html?id=79555"
    I = "//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg"
  }
  3 = @{
    B = "2024.03.09"
    C = "景德镇·江报国风动漫展 "
    D = "迎宾大道与寺山路交叉口东200米 陶博城"
    E = "2024.03.09 09:00-03.10 17:00"
    F = 956
    G = 55
    H = "https://show.bilibili.com/platform/detail.html?id=81362"
    I = "//i2.hdslb.com/bfs/openplatform/202402/oM49o66R1708334630235.jpeg"
  }
  4 = @{
    B = "2024.03.16"
    C = "景德镇·原神X崩铁X崩坏动漫展only"
    D = "陶阳南路188号 晨枫臻品酒店"
    E = "2024.03.16 10:00-03.16 17:00"
    F = 60
    G = 55
    H = "https://show.bilibili.com/platform/detail.html?id=80920"
    I = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"
  }
  5 = @{
    B = "2024.03.16"
    C = "江西·ShiningStaR动漫游戏文化节5th"
    D = "高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆"
    E = "2024.03.16 09:30-03.17 17:00"
    F = 2186
    G = 60
    H = "https://show.bilibili.com/platform/detail.html?id=81792"
    I = "//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg"
  }
  6 = @{
    B = "2024.03.23"
    C = "上饶·原×铁×崩only"
    D = "五三东大道42号 回禾酒店"
    E = "2024.03.23 10:00-03.23 17:00"
    F = 37
    G = 60
    H = "https://show.bilibili.com/platform/detail.html?id=81103"
    I = "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"
  }
  7 = @{
    B = "2024.03.23"
    C = "南昌·AP动漫游戏嘉年华"
    D = "八一桥街道青山南路118号 蓝海会展中心"
    E = "2024.03.23 09:00-03.24 17:00"
    F = 1337
    G = 58.5
    H = "https://show.bilibili.com/platform/detail.html?id=81232"
    I = "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"
  }
  8 = @{
    B = "2024.03.23"
    C = "南昌·原X穹X崩only"
    D = "丰和北大道299号 新吉花园酒店"
    E = "2024.03.23 10:00-03.23 17:00"
    F = 63
    G = 65
    H = "https://show.bilibili.com/platform/detail.html?id=80807"
    I = "//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg"
  }
  9 = @{
    B = "2024.03.23"
    C = "南昌·运动番only春季集训"
    D = "创新三路777号 南昌小飞侠章鱼文化体育公园"
    E = "2024.03.23 10:00-03.24 17:00"
    F = 135
    G = 58
    H = "https://show.bilibili.com/platform/detail.html?id=81950"
    I = "//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg"
  }
  10 = @{
    B = "2024.03.24"
    C = "南昌·AP动漫游戏  嘉年华内场票-小N&子音"
    D = "八一桥街道青山南路118号 蓝海会展中心"
    E = "2024.03.24 09:00-03.24 17:00"
    F = 43
    G = 218
    H = "https://show.bilibili.com/platform/detail.html?id=81973"
    I = "//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg"
  }
  11 = @{
    B = "2024.03.30"
    C = "南昌·CM01动漫游戏博览会"
    D = "怀玉山大道1315号 南昌绿地国际博览中心"
    E = "2024.03.30 10:00-03.31 17:00"
    F = 321
    G = 55
    H = "https://show.bilibili.com/platform/detail.html?id=81691"
    I = "//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png"
  }
  12 = @{
    B = "2024.03.30"
    C = "鹰潭·原×铁×崩only"
    D = "南站路24号 回禾酒店(鹰潭火车站店)"
    E = "2024.03.30 10:00-03.30 17:00"
    F = 25
    G = 60
    H = "https://show.bilibili.com/platform/detail.html?id=81097"
    I = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"
  }
}

$targetSheetNames = @("展览", "全部类型")

foreach ($sheetName in $targetSheetNames) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Column B holds plain-text "YYYY.MM.DD" strings. Force text format
  # first so Excel does not reinterpret them as date serials.
  $ws.Range("B2:B12").NumberFormat = "@"

  foreach ($rowNum in 2..12) {
    $data = $rowsAfter[$rowNum]
    $ws.Cells.Item($rowNum, 2).Value = $data.B
    $ws.Cells.Item($rowNum, 3).Value = $data.C
    $ws.Cells.Item($rowNum, 4).Value = $data.D
    $ws.Cells.Item($rowNum, 5).Value = $data.E
    $ws.Cells.Item($rowNum, 6).Value = $data.F
    $ws.Cells.Item($rowNum, 7).Value = $data.G
    $ws.Cells.Item($rowNum, 8).Value = $data.H
    $ws.Cells.Item($rowNum, 9).Value = $data.I
  }

  # The old row 13 ("鹰潭·原×铁×崩only" before the shift) is now a
  # duplicate of row 12 content-wise; delete it so the table is back
  # to 12 data rows (dimension A1:I12).
  $ws.Rows("13").Delete()
}
